$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'285.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'2.34%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.90%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.051"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.45%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06684"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.69%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.340"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.25%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'2.40%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.351"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'3.88%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9388"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'5.13%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1577"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.53%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.06630"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'13.62%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07657"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.00%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02915"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.08982"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.01%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001572"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.56%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04465"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.83%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0006458"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.21%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006559"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'6.47%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.484"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.33%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.232"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-3.17%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.3208"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.93%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1307"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-3.08%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.045"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'3.73%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1523"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.19%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001178"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.19%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004480"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001244"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'5.53%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001613"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'-2.24%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04202"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.89%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006723"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.10%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1248"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-11.75%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002011"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.32%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01219"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'9.24%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005671"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.38%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'20.74%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01303"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-29.39%"
$ws.Range("E47").Style = "Normal"
